# feat: add 2022-Q1 data
#
# The workbook tracks quarterly fund-holdings snapshots, one worksheet per
# quarter, plus a trailing "总计" (totals) summary sheet. This adds a new
# "2022-Q1" quarter sheet (positioned right before "总计") and appends a
# fresh "总计" sheet after it with an extra leading row for the new quarter.

$wb = $excel.ActiveWorkbook

# --- 1. Re-purpose the current "总计" sheet as the new "2022-Q1" sheet
#        (keeps its identity/sheetId), and append a brand-new sheet right
#        after it to become the new "总计".
$quarterSheet = $wb.Worksheets.Item("总计")
$quarterSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add($null, $quarterSheet)
$totalSheet.Name = "总计"

# --- 2. Snapshot the old "总计" table (still living in what is now named
#        "2022-Q1") into the new "总计" sheet before it gets overwritten:
#        header row stays at row 1, the data rows shift down one row to
#        make room for the new 2022-Q1 entry, and the running A-column
#        index is bumped accordingly.
$quarterSheet.Range("B1:D1").Copy($totalSheet.Range("B1"))
$quarterSheet.Range("A2").Copy($totalSheet.Range("A2"))
$quarterSheet.Range("A2:D6").Copy($totalSheet.Range("A3"))
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

# --- 3. Build the "2022-Q1" fund-holding sheet, matching the layout of
#        the other quarter sheets (e.g. "2021-Q4"): pull header + row
#        formatting from that template, then fill in this quarter's data.
$quarterSheet.Cells.Clear()
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H2").Copy($quarterSheet.Range("A1"))

$textCells = $quarterSheet.Range("B2:G2")
$textCells.NumberFormat = "@"
$quarterSheet.Range("B2").Value = "162416"
$quarterSheet.Range("C2").Value = "华宝港股通恒生香港35指数(LOF)"
$quarterSheet.Range("D2").Value = "0.21"
$quarterSheet.Range("E2").Value = "94.50"
$quarterSheet.Range("F2").Value = "3.75"
$quarterSheet.Range("G2").Value = "0.0079"
$quarterSheet.Range("H2").Value = 9
